# Generate Report for Handoff
# A new source file (ffff563dc99b-6133-41f8-8051-44718e383461.md) was handed off
# together with the existing 3b7b5ac1-8d6e-4fdc-a32c-a5af63904b3a.md (which itself
# is a renamed/re-handed-off version of the old 257eb94a-064d-4a6e-b247-8783393e957e.md).
# This adds a new row to each of the three worksheets (Overview, zh-cn, de-de) and
# refreshes the handoff file name / handoff datetime for the existing row.

$wb = $excel.ActiveWorkbook

$mdBase     = "https://github.com/OpenLocalizationTest/oltest/blob/47da51ae6e4fa9a4be762d3c2beeb0467bd03593/e2e/"
$configUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/47da51ae6e4fa9a4be762d3c2beeb0467bd03593/.localization-config"
$zhcnXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9e7af1bbe963e77072c81d9913db7c84f81b4d5a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/"
$dedeXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d4c5488fd57d8b3b5eea5049e73f38156a0b5aa9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/"

$newMd1 = "3b7b5ac1-8d6e-4fdc-a32c-a5af63904b3a.md"
$newMd2 = "ffff563dc99b-6133-41f8-8051-44718e383461.md"
$configName = ".localization-config"

$zhcnXlf = "3b7b5ac1-8d6e-4fdc-a32c-a5af63904b3a.ecb45406ec952809a85e2b0699a83d4f8badbcd9.zh-cn.xlf"
$dedeXlf = "3b7b5ac1-8d6e-4fdc-a32c-a5af63904b3a.ecb45406ec952809a85e2b0699a83d4f8badbcd9.de-de.xlf"

$zhcnHandoffDatetime = "2016-03-07 04:59:19"
$dedeHandoffDatetime = "2016-03-07 04:59:30"
$epoch = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A1:C10").Hyperlinks.Delete()

$ws1.Range("A2").Value = $newMd1
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

$ws1.Range("A3").Value = $newMd2
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Range("A4").Value = $configName
$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

$ws1.Hyperlinks.Add($ws1.Range("A2"), ($mdBase + $newMd1), "", "", $newMd1)
$ws1.Hyperlinks.Add($ws1.Range("A3"), ($mdBase + $newMd2), "", "", $newMd2)
$ws1.Hyperlinks.Add($ws1.Range("A4"), $configUrl, "", "", $configName)

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A1:I10").Hyperlinks.Delete()

$ws2.Range("A2").Value = $newMd1
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("C2").Value = $zhcnXlf
$ws2.Range("D2").Value = $zhcnHandoffDatetime
$ws2.Range("G2").Value = $epoch
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = $newMd2
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("C3").Value = $zhcnXlf
$ws2.Range("D3").Value = $zhcnHandoffDatetime
$ws2.Range("G3").Value = $epoch
$ws2.Range("H3").Value = "Include"

$ws2.Range("A4").Value = $configName
$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = $epoch
$ws2.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G4").Value = $epoch
$ws2.Range("H4").Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Range("A2"), ($mdBase + $newMd1), "", "", $newMd1)
$ws2.Hyperlinks.Add($ws2.Range("C2"), ($zhcnXlfBase + $zhcnXlf), "", "", $zhcnXlf)
$ws2.Hyperlinks.Add($ws2.Range("A3"), ($mdBase + $newMd2), "", "", $newMd2)
$ws2.Hyperlinks.Add($ws2.Range("C3"), ($zhcnXlfBase + $zhcnXlf), "", "", $zhcnXlf)
$ws2.Hyperlinks.Add($ws2.Range("A4"), $configUrl, "", "", $configName)

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A1:I10").Hyperlinks.Delete()

$ws3.Range("A2").Value = $newMd1
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("C2").Value = $dedeXlf
$ws3.Range("D2").Value = $dedeHandoffDatetime
$ws3.Range("G2").Value = $epoch
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = $newMd2
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("C3").Value = $dedeXlf
$ws3.Range("D3").Value = $dedeHandoffDatetime
$ws3.Range("G3").Value = $epoch
$ws3.Range("H3").Value = "Include"

$ws3.Range("A4").Value = $configName
$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = $epoch
$ws3.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G4").Value = $epoch
$ws3.Range("H4").Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Range("A2"), ($mdBase + $newMd1), "", "", $newMd1)
$ws3.Hyperlinks.Add($ws3.Range("C2"), ($dedeXlfBase + $dedeXlf), "", "", $dedeXlf)
$ws3.Hyperlinks.Add($ws3.Range("A3"), ($mdBase + $newMd2), "", "", $newMd2)
$ws3.Hyperlinks.Add($ws3.Range("C3"), ($dedeXlfBase + $dedeXlf), "", "", $dedeXlf)
$ws3.Hyperlinks.Add($ws3.Range("A4"), $configUrl, "", "", $configName)
